# Fix port numbers in slide deck.
#
# Slide 8 ("Labels & Named Ports") shows two Service "cards":
#   - shape id 7  -> the "client" service card (TargetPort should be 8080)
#   - shape id 27 -> the "nginx" service card (Port should be 80)
# Both had the port numbers swapped; this script corrects the text and
# also restores the original shape (z-order) ordering seen in the fixed
# deck: Rect(7), Rect(27), Title(2), Group(31), Group(32), Connector(16),
# Connector(28).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(8)

function Get-ShapeById($slide, $id) {
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $candidate = $slide.Shapes.Item($i)
        if ($candidate.Id -eq $id) {
            return $candidate
        }
    }
    return $null
}

# --- Text fixes -----------------------------------------------------------

# Shape 7 ("client" service card): "TargetPort: 80" -> "TargetPort: 8080"
$sh7 = Get-ShapeById $s 7
$para3 = $sh7.TextFrame.TextRange.Paragraphs(3)
$para3.Runs(2).Text = "8080"

# Shape 27 ("nginx" service card): "Port: 8080" -> "Port: 80"
$sh27 = Get-ShapeById $s 27
$para5 = $sh27.TextFrame.TextRange.Paragraphs(5)
$para5.Text = "Port: 80"

# --- Shape order fix --------------------------------------------------------
# Desired back-to-front order: 7, 27, 2, 31, 32, 16, 28
# Send to back in reverse of the desired order so each one lands correctly.
$desiredBackToFront = @(28, 16, 32, 31, 2, 27, 7)
foreach ($id in $desiredBackToFront) {
    $shape = Get-ShapeById $s $id
    $shape.ZOrder(1)  # msoSendToBack
}
